$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 3-39 (columns B:G).
# topic | titolo | giornale | social | negativo | positivo
$data = @(
    @('POLITICA','Attacco di Hamas a Israele','FanPage','Facebook',62,38),
    @('POLITICA','Attacco di Hamas a Israele','FanPage','Instagram',84,16),
    @('POLITICA','Attacco di Hamas a Israele','FanPage','YouTube',83,17),
    @('POLITICA','Attacco di Hamas a Israele','Il Corriere Della Sera','Facebook',82,18),
    @('POLITICA','Attacco di Hamas a Israele','Il Corriere Della Sera','Instagram',84,16),
    @('POLITICA','Attacco di Hamas a Israele','Il Corriere Della Sera','YouTube',81,19),
    @('POLITICA','Attacco di Hamas a Israele','La Repubblica','Facebook',86,13),
    @('POLITICA','Attacco di Hamas a Israele','La Repubblica','Instagram',85,8),
    @('POLITICA','Attacco di Hamas a Israele','La Repubblica','YouTube',78,17),
    @('POLITICA','Invasione Ucraina','FanPage','Facebook',94,6),
    @('POLITICA','Invasione Ucraina','FanPage','Instagram',91,9),
    @('POLITICA','Invasione Ucraina','FanPage','YouTube',79,21),
    @('POLITICA','Invasione Ucraina','Il Corriere Della Sera','Facebook',87,13),
    @('POLITICA','Invasione Ucraina','Il Corriere Della Sera','Instagram',89,11),
    @('POLITICA','Invasione Ucraina','Il Corriere Della Sera','YouTube',80,20),
    @('POLITICA','Invasione Ucraina','La Repubblica','Facebook',60,40),
    @('POLITICA','Invasione Ucraina','La Repubblica','Instagram',84,16),
    @('POLITICA','Invasione Ucraina','La Repubblica','YouTube',78,22),
    @('POLITICA','L''incoronazione di Re Carlo','FanPage','Facebook',45,55),
    @('POLITICA','L''incoronazione di Re Carlo','FanPage','Instagram',62,38),
    @('POLITICA','L''incoronazione di Re Carlo','FanPage','YouTube',67,33),
    @('POLITICA','L''incoronazione di Re Carlo','Il Corriere Della Sera','Facebook',69,31),
    @('POLITICA','L''incoronazione di Re Carlo','Il Corriere Della Sera','Instagram',78,22),
    @('POLITICA','L''incoronazione di Re Carlo','Il Corriere Della Sera','YouTube',75,25),
    @('POLITICA','L''incoronazione di Re Carlo','La Repubblica','Facebook',65,35),
    @('POLITICA','L''incoronazione di Re Carlo','La Repubblica','Instagram',72,28),
    @('POLITICA','L''incoronazione di Re Carlo','La Repubblica','YouTube',66,34),
    @('POLITICA','La morte di Silvio Berlusconi','FanPage','Facebook',74,26),
    @('POLITICA','La morte di Silvio Berlusconi','FanPage','Instagram',75,25),
    @('POLITICA','La morte di Silvio Berlusconi','FanPage','YouTube',77,23),
    @('POLITICA','La morte di Silvio Berlusconi','Il Corriere Della Sera','Facebook',75,24),
    @('POLITICA','La morte di Silvio Berlusconi','Il Corriere Della Sera','Instagram',82,18),
    @('POLITICA','La morte di Silvio Berlusconi','Il Corriere Della Sera','YouTube',65,35),
    @('POLITICA','La morte di Silvio Berlusconi','Il Corrirere Della Sera','Facebook',1,0),
    @('POLITICA','La morte di Silvio Berlusconi','La Repubblica','Facebook',74,26),
    @('POLITICA','La morte di Silvio Berlusconi','La Repubblica','Instagram',81,19),
    @('POLITICA','La morte di Silvio Berlusconi','La Repubblica','YouTube',66,34)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = ""
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
}
